$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56
$ws.Range("B56").Value = "'33.59"
$ws.Range("C56").Value = "'0.39"
$ws.Range("D56").Value = "'5.12"
$ws.Range("E56").Value = "'1.47"

# Row 57
$ws.Range("B57").Value = "'14.66"
$ws.Range("C57").Value = "'0.49"
$ws.Range("D57").Value = "'4.46"
$ws.Range("E57").Value = "'1.49"

# Row 58
$ws.Range("B58").Value = "'8.24"
$ws.Range("C58").Value = "'0.44"
$ws.Range("D58").Value = "'3.24"
$ws.Range("E58").Value = "'1.40"

# Row 59
$ws.Range("B59").Value = "'3.60"
$ws.Range("C59").Value = "'0.47"
$ws.Range("D59").Value = "'3.26"
$ws.Range("E59").Value = "'1.42"

# Row 60
$ws.Range("B60").Value = "'4.75"
$ws.Range("C60").Value = "'0.55"
$ws.Range("D60").Value = "'3.68"
$ws.Range("E60").Value = "'1.52"

# Row 61
$ws.Range("B61").Value = "'4.87"
$ws.Range("C61").Value = "'0.43"
$ws.Range("D61").Value = "'5.12"
$ws.Range("E61").Value = "'1.73"
